$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liste")

$ws.Range("F3").Value = "U3-4"
$ws.Range("F4").Value = "U3-Amphi"
$ws.Range("F7").Value = "U3-4"
$ws.Range("F8").Value = "U3-4"
$ws.Range("F16").Value = "U3-Amphi"
